$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.220.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.399.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.398.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.636"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.927.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.383.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.040.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "89.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "580.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.49%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.379"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0756"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.126.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0420"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.134"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.48%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.20%  "
